$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 05:55:25"
$ws1.Cells.Item(3, 1).Value = "Total filas: 44"

$ws1.Cells.Item(18, 1).Value = "05:55:25"
$ws1.Cells.Item(18, 2).Value = "05:55"
$ws1.Cells.Item(18, 3).Value = "10_OLMOS"
$ws1.Cells.Item(18, 4).Value = 0
$ws1.Cells.Item(18, 5).Value = "LP1912"

$ws1.Cells.Item(19, 1).Value = "05:55:25"
$ws1.Cells.Item(19, 2).Value = "05:56"
$ws1.Cells.Item(19, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(19, 4).Value = 1
$ws1.Cells.Item(19, 5).Value = "LP1912"

$ws1.Cells.Item(20, 1).Value = "05:22:24"
$ws1.Cells.Item(20, 2).Value = "06:04"
$ws1.Cells.Item(20, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(20, 4).Value = 42
$ws1.Cells.Item(20, 5).Value = "LP1912"

$ws1.Cells.Item(21, 1).Value = "04:18:53"
$ws1.Cells.Item(21, 2).Value = "06:05"
$ws1.Cells.Item(21, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(21, 4).Value = 107
$ws1.Cells.Item(21, 5).Value = "LP1912"

$ws1.Cells.Item(22, 1).Value = "05:55:25"
$ws1.Cells.Item(22, 2).Value = "06:11"
$ws1.Cells.Item(22, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(22, 4).Value = 16
$ws1.Cells.Item(22, 5).Value = "LP1912"

$ws1.Cells.Item(23, 1).Value = "04:56:06"
$ws1.Cells.Item(23, 2).Value = "06:12"
$ws1.Cells.Item(23, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(23, 4).Value = 76
$ws1.Cells.Item(23, 5).Value = "LP1912"

$ws1.Cells.Item(24, 1).Value = "05:55:25"
$ws1.Cells.Item(24, 2).Value = "06:13"
$ws1.Cells.Item(24, 3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(24, 4).Value = 91
$ws1.Cells.Item(24, 5).Value = "LP1912"

$ws1.Cells.Item(25, 1).Value = "05:55:25"
$ws1.Cells.Item(25, 2).Value = "06:14"
$ws1.Cells.Item(25, 3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(25, 4).Value = 19
$ws1.Cells.Item(25, 5).Value = "LP1912"

$ws1.Cells.Item(26, 1).Value = "04:42:52"
$ws1.Cells.Item(26, 2).Value = "06:20"
$ws1.Cells.Item(26, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(26, 4).Value = 98
$ws1.Cells.Item(26, 5).Value = "LP1912"

$ws1.Cells.Item(27, 1).Value = "05:55:25"
$ws1.Cells.Item(27, 2).Value = "06:21"
$ws1.Cells.Item(27, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(27, 4).Value = 26
$ws1.Cells.Item(27, 5).Value = "LP1912"

$ws1.Cells.Item(28, 1).Value = "04:42:52"
$ws1.Cells.Item(28, 2).Value = "06:26"
$ws1.Cells.Item(28, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(28, 4).Value = 104
$ws1.Cells.Item(28, 5).Value = "LP1912"

$ws1.Cells.Item(29, 1).Value = "05:55:25"
$ws1.Cells.Item(29, 2).Value = "06:27"
$ws1.Cells.Item(29, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(29, 4).Value = 32
$ws1.Cells.Item(29, 5).Value = "LP1912"

$ws1.Cells.Item(30, 1).Value = "05:55:25"
$ws1.Cells.Item(30, 2).Value = "06:29"
$ws1.Cells.Item(30, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(30, 4).Value = 34
$ws1.Cells.Item(30, 5).Value = "LP1912"

$ws1.Cells.Item(31, 1).Value = "04:56:06"
$ws1.Cells.Item(31, 2).Value = "06:30"
$ws1.Cells.Item(31, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(31, 4).Value = 94
$ws1.Cells.Item(31, 5).Value = "LP1912"

$ws1.Cells.Item(32, 1).Value = "05:55:25"
$ws1.Cells.Item(32, 2).Value = "06:31"
$ws1.Cells.Item(32, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(32, 4).Value = 36
$ws1.Cells.Item(32, 5).Value = "LP1912"

$ws1.Cells.Item(33, 1).Value = "05:55:25"
$ws1.Cells.Item(33, 2).Value = "06:44"
$ws1.Cells.Item(33, 3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(33, 4).Value = 49
$ws1.Cells.Item(33, 5).Value = "LP1912"

$ws1.Cells.Item(34, 1).Value = "05:55:25"
$ws1.Cells.Item(34, 2).Value = "06:46"
$ws1.Cells.Item(34, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(34, 4).Value = 51
$ws1.Cells.Item(34, 5).Value = "LP1912"

$ws1.Cells.Item(35, 1).Value = "04:56:06"
$ws1.Cells.Item(35, 2).Value = "06:47"
$ws1.Cells.Item(35, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(35, 4).Value = 111
$ws1.Cells.Item(35, 5).Value = "LP1912"

$ws1.Cells.Item(36, 1).Value = "05:55:25"
$ws1.Cells.Item(36, 2).Value = "06:59"
$ws1.Cells.Item(36, 3).Value = "14_ABASTO"
$ws1.Cells.Item(36, 4).Value = 64
$ws1.Cells.Item(36, 5).Value = "LP1912"

$ws1.Cells.Item(37, 1).Value = "05:55:25"
$ws1.Cells.Item(37, 2).Value = "07:04"
$ws1.Cells.Item(37, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(37, 4).Value = 69
$ws1.Cells.Item(37, 5).Value = "LP1912"

$ws1.Cells.Item(38, 1).Value = "05:55:25"
$ws1.Cells.Item(38, 2).Value = "07:05"
$ws1.Cells.Item(38, 3).Value = "15_ABASTO"
$ws1.Cells.Item(38, 4).Value = 70
$ws1.Cells.Item(38, 5).Value = "LP1912"

$ws1.Cells.Item(39, 1).Value = "05:55:25"
$ws1.Cells.Item(39, 2).Value = "07:07"
$ws1.Cells.Item(39, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(39, 4).Value = 72
$ws1.Cells.Item(39, 5).Value = "LP1912"

$ws1.Cells.Item(40, 1).Value = "05:55:25"
$ws1.Cells.Item(40, 2).Value = "07:11"
$ws1.Cells.Item(40, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(40, 4).Value = 76
$ws1.Cells.Item(40, 5).Value = "LP1912"

$ws1.Cells.Item(41, 1).Value = "05:55:25"
$ws1.Cells.Item(41, 2).Value = "07:15"
$ws1.Cells.Item(41, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(41, 4).Value = 80
$ws1.Cells.Item(41, 5).Value = "LP1912"

$ws1.Cells.Item(42, 1).Value = "05:55:25"
$ws1.Cells.Item(42, 2).Value = "07:21"
$ws1.Cells.Item(42, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(42, 4).Value = 86
$ws1.Cells.Item(42, 5).Value = "LP1912"

$ws1.Cells.Item(43, 1).Value = "05:55:25"
$ws1.Cells.Item(43, 2).Value = "07:31"
$ws1.Cells.Item(43, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(43, 4).Value = 96
$ws1.Cells.Item(43, 5).Value = "LP1912"

$ws1.Cells.Item(44, 1).Value = "05:55:25"
$ws1.Cells.Item(44, 2).Value = "07:31"
$ws1.Cells.Item(44, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(44, 4).Value = 96
$ws1.Cells.Item(44, 5).Value = "LP1912"

$ws1.Cells.Item(45, 1).Value = "05:55:25"
$ws1.Cells.Item(45, 2).Value = "07:32"
$ws1.Cells.Item(45, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(45, 4).Value = 97
$ws1.Cells.Item(45, 5).Value = "LP1912"

$ws1.Cells.Item(46, 1).Value = "05:55:25"
$ws1.Cells.Item(46, 2).Value = "07:36"
$ws1.Cells.Item(46, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(46, 4).Value = 101
$ws1.Cells.Item(46, 5).Value = "LP1912"

$ws1.Cells.Item(47, 1).Value = "05:55:25"
$ws1.Cells.Item(47, 2).Value = "07:39"
$ws1.Cells.Item(47, 3).Value = "10_OLMOS"
$ws1.Cells.Item(47, 4).Value = 104
$ws1.Cells.Item(47, 5).Value = "LP1912"

$ws1.Cells.Item(48, 1).Value = "05:55:25"
$ws1.Cells.Item(48, 2).Value = "07:47"
$ws1.Cells.Item(48, 3).Value = "14_ABASTO"
$ws1.Cells.Item(48, 4).Value = 112
$ws1.Cells.Item(48, 5).Value = "LP1912"

$ws1.Cells.Item(49, 1).Value = "05:55:25"
$ws1.Cells.Item(49, 2).Value = "07:51"
$ws1.Cells.Item(49, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(49, 4).Value = 116
$ws1.Cells.Item(49, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 05:55:25"
$ws2.Cells.Item(3, 1).Value = "Total filas: 10"
$ws2.Cells.Item(10, 1).Value = "05:55:25"
$ws2.Cells.Item(10, 4).Value = 16
$ws2.Cells.Item(12, 1).Value = "05:55:25"
$ws2.Cells.Item(12, 4).Value = 51
$ws2.Cells.Item(14, 1).Value = "05:55:25"
$ws2.Cells.Item(14, 4).Value = 76
$ws2.Cells.Item(15, 1).Value = "05:55:25"
$ws2.Cells.Item(15, 2).Value = "07:51"
$ws2.Cells.Item(15, 3).Value = "215D_EL PATO"
$ws2.Cells.Item(15, 4).Value = 116
$ws2.Cells.Item(15, 5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 05:55:25"
$ws3.Cells.Item(3, 1).Value = "Total filas: 8"
$ws3.Cells.Item(9, 1).Value = "05:55:25"
$ws3.Cells.Item(9, 4).Value = 14
$ws3.Cells.Item(11, 1).Value = "05:55:25"
$ws3.Cells.Item(11, 4).Value = 38
$ws3.Cells.Item(12, 1).Value = "05:55:25"
$ws3.Cells.Item(12, 4).Value = 65
$ws3.Cells.Item(13, 1).Value = "05:55:25"
$ws3.Cells.Item(13, 2).Value = "07:35"
$ws3.Cells.Item(13, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(13, 4).Value = 100
$ws3.Cells.Item(13, 5).Value = "L6173"

Write-Host "Edit complete"